{"js": "// Lattice multiplication exercise refresh: replace the 5 lines of text\n// (problem, multiplicand digits, separator, and two lattice row labels)\n// inside every cell of the 5x3 exercise table with newly generated values,\n// while preserving the existing paragraph/run/formatting structure\n// (single run per cell, sz=32, <w:t> segments joined by <w:br/>).\n\nconst CELL_LINES = [[\"67 x 38\", \"  3    8\", \"  ----\", \"6|    |\", \"7|    |\"], [\"66 x 88\", \"  8    8\", \"  ----\", \"6|    |\", \"6|    |\"], [\"28 x 97\", \"  9    7\", \"  ----\", \"2|    |\", \"8|    |\"], [\"73 x 28\", \"  2    8\", \"  ----\", \"7|    |\", \"3|    |\"], [\"53 x 50\", \"  5    0\", \"  ----\", \"5|    |\", \"3|    |\"], [\"53 x 86\", \"  8    6\", \"  ----\", \"5|    |\", \"3|    |\"], [\"96 x 53\", \"  5    3\", \"  ----\", \"9|    |\", \"6|    |\"], [\"45 x 20\", \"  2    0\", \"  ----\", \"4|    |\", \"5|    |\"], [\"45 x 34\", \"  3    4\", \"  ----\", \"4|    |\", \"5|    |\"], [\"50 x 98\", \"  9    8\", \"  ----\", \"5|    |\", \"0|    |\"], [\"55 x 55\", \"  5    5\", \"  ----\", \"5|    |\", \"5|    |\"], [\"25 x 69\", \"  6    9\", \"  ----\", \"2|    |\", \"5|    |\"], [\"54 x 45\", \"  4    5\", \"  ----\", \"5|    |\", \"4|    |\"], [\"29 x 14\", \"  1    4\", \"  ----\", \"2|    |\", \"9|    |\"], [\"13 x 44\", \"  4    4\", \"  ----\", \"1|    |\", \"3|    |\"]];\n\nfunction escapeXml(s) {\n  return s.replace(/&/g, \"&amp;\")\n          .replace(/</g, \"&lt;\")\n          .replace(/>/g, \"&gt;\")\n          .replace(/\"/g, \"&quot;\")\n          .replace(/'/g, \"&apos;\");\n}\n\nfunction cellOoxml(lines) {\n  const runsXml = lines\n    .map((line) => `<w:t xml:space=\"preserve\">${escapeXml(line)}</w:t>`)\n    .join(\"<w:br/>\");\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n<pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n</Relationships>\n</pkg:xmlData>\n</pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:r><w:rPr><w:sz w:val=\"32\"/></w:rPr>${runsXml}</w:r></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst columnCount = CELL_LINES.length / table.rowCount;\n\nfor (let i = 0; i < CELL_LINES.length; i++) {\n  const rowIndex = Math.floor(i / columnCount);\n  const colIndex = i % columnCount;\n  const cell = table.getCell(rowIndex, colIndex);\n  cell.body.insertOoxml(cellOoxml(CELL_LINES[i]), \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Lattice multiplication exercise refresh: replace the 5 lines of text\n# (problem, multiplicand digits, separator, and two lattice row labels)\n# inside every cell of the 5x3 exercise table with newly generated values.\n# Setting Cell.Range.Text with embedded vertical-tab (chr(11)) separators\n# reproduces Word's own <w:br/>-separated single-run cell layout and keeps\n# the existing run formatting (sz=32).\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$cellLines = @(\n  @('67 x 38', '  3    8', '  ----', '6|    |', '7|    |'),\n  @('66 x 88', '  8    8', '  ----', '6|    |', '6|    |'),\n  @('28 x 97', '  9    7', '  ----', '2|    |', '8|    |'),\n  @('73 x 28', '  2    8', '  ----', '7|    |', '3|    |'),\n  @('53 x 50', '  5    0', '  ----', '5|    |', '3|    |'),\n  @('53 x 86', '  8    6', '  ----', '5|    |', '3|    |'),\n  @('96 x 53', '  5    3', '  ----', '9|    |', '6|    |'),\n  @('45 x 20', '  2    0', '  ----', '4|    |', '5|    |'),\n  @('45 x 34', '  3    4', '  ----', '4|    |', '5|    |'),\n  @('50 x 98', '  9    8', '  ----', '5|    |', '0|    |'),\n  @('55 x 55', '  5    5', '  ----', '5|    |', '5|    |'),\n  @('25 x 69', '  6    9', '  ----', '2|    |', '5|    |'),\n  @('54 x 45', '  4    5', '  ----', '5|    |', '4|    |'),\n  @('29 x 14', '  1    4', '  ----', '2|    |', '9|    |'),\n  @('13 x 44', '  4    4', '  ----', '1|    |', '3|    |')\n)\n\n$lineBreak = [char]11\n$columnCount = $table.Columns.Count\n\nfor ($i = 0; $i -lt $cellLines.Count; $i++) {\n  $rowIndex = [int][Math]::Floor($i / $columnCount) + 1\n  $colIndex = ($i % $columnCount) + 1\n  $cell = $table.Cell($rowIndex, $colIndex)\n  $cell.Range.Text = [string]::Join($lineBreak, $cellLines[$i])\n}\n\n"}
